# Update countries & provincias Spain
# - Refresh the "last updated" timestamp
# - Refresh the COVID-19 country stats (which re-sorts several rows by
#   total cases, so the country names in column A of the affected rows
#   change along with their B:H figures)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp footer (A1)
$ws.Range("A1").Value = "Datos actualizados a 10 de Mayo de 2020 a las 09:34"

# Rows 28-29: Singapur overtakes Irlanda
$ws.Range("A28").Value = "Singapur"
$ws.Range("B28").Value = 23336
$ws.Range("C28").Value = 876
$ws.Range("D28").Value = 2296
$ws.Range("E28").Value = 21020
$ws.Range("F28").Value = 23
$ws.Range("G28").Value = 0
$ws.Range("H28").Value = 20

$ws.Range("A29").Value = "Irlanda"
$ws.Range("B29").Value = 22760
$ws.Range("C29").Value = 0
$ws.Range("D29").Value = 17110
$ws.Range("E29").Value = 4204
$ws.Range("F29").Value = 72
$ws.Range("G29").Value = 0
$ws.Range("H29").Value = 1446

# Row 36 (Polonia): updated recuperados / casos activos, name unchanged
$ws.Range("D36").Value = 5698
$ws.Range("E36").Value = 9168

# Rows 37-38: Ucrania overtakes Rumania
$ws.Range("A37").Value = "Ucrania"
$ws.Range("B37").Value = 15232
$ws.Range("C37").Value = 522
$ws.Range("D37").Value = 3060
$ws.Range("E37").Value = 11781
$ws.Range("F37").Value = 201
$ws.Range("G37").Value = 15
$ws.Range("H37").Value = 391

$ws.Range("A38").Value = "Rumania"
$ws.Range("B38").Value = 15131
$ws.Range("C38").Value = 0
$ws.Range("D38").Value = 6912
$ws.Range("E38").Value = 7280
$ws.Range("F38").Value = 245
$ws.Range("G38").Value = 0
$ws.Range("H38").Value = 939

# Rows 66-68: Armenia overtakes Hungria and Oman
$ws.Range("A66").Value = "Armenia"
$ws.Range("B66").Value = 3313
$ws.Range("C66").Value = 138
$ws.Range("D66").Value = 1325
$ws.Range("E66").Value = 1943
$ws.Range("F66").Value = 10
$ws.Range("G66").Value = 1
$ws.Range("H66").Value = 45

$ws.Range("A67").Value = "Hungria"
$ws.Range("B67").Value = 3263
$ws.Range("C67").Value = 50
$ws.Range("D67").Value = 933
$ws.Range("E67").Value = 1917
$ws.Range("F67").Value = 50
$ws.Range("G67").Value = 8
$ws.Range("H67").Value = 413

$ws.Range("A68").Value = "Oman"
$ws.Range("B68").Value = 3224
$ws.Range("C68").Value = 0
$ws.Range("D68").Value = 1068
$ws.Range("E68").Value = 2139
$ws.Range("F68").Value = 17
$ws.Range("G68").Value = 0
$ws.Range("H68").Value = 17

# Rows 99-100: Letonia overtakes Consejo Danes para los Refugiados
$ws.Range("A99").Value = "Letonia"
$ws.Range("B99").Value = 939
$ws.Range("C99").Value = 9
$ws.Range("D99").Value = 464
$ws.Range("E99").Value = 457
$ws.Range("F99").Value = 2
$ws.Range("G99").Value = 0
$ws.Range("H99").Value = 18

$ws.Range("A100").Value = "Consejo Danes para los Refugiados"
$ws.Range("B100").Value = 937
$ws.Range("C100").Value = 0
$ws.Range("D100").Value = 130
$ws.Range("E100").Value = 768
$ws.Range("F100").Value = 0
$ws.Range("G100").Value = 0
$ws.Range("H100").Value = 39

# Row 152 (Suazilandia): updated recuperados / casos activos, name unchanged
$ws.Range("D152").Value = 143
$ws.Range("E152").Value = 3

# Rows 206-207: Montserrat overtakes Seychelles
$ws.Range("A206").Value = "Montserrat"
$ws.Range("D206").Value = 8
$ws.Range("E206").Value = 2
$ws.Range("F206").Value = 1
$ws.Range("H206").Value = 1

$ws.Range("A207").Value = "Seychelles"
$ws.Range("D207").Value = 10
$ws.Range("E207").Value = 1
$ws.Range("F207").Value = 0
$ws.Range("H207").Value = 0
